# Add "2022-Q1" sheet: the detailed fund-holdings table that used to be
# on the "总计" sheet moves to a brand-new "2022-Q1" sheet (inserted
# right before "总计"), and "总计" is rebuilt as the quarterly summary
# table with an additional first row for 2022-Q1.
#
# NOTE: worksheet references captured *before* a sheet-collection
# structural change (Add/Copy/Delete/rename) can end up pointing at the
# wrong sheet afterwards because they are resolved by position. To stay
# safe, sheets are re-fetched by name after every structural change.

$wb = $excel.ActiveWorkbook

$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")

# Create the new "2022-Q1" worksheet by duplicating "2021-Q4" (so it
# inherits the same sheet-level formatting/margins/outline settings)
# and placing it immediately before "总计".
$q4Sheet.Copy($totalSheet)

$newQ1 = $wb.Worksheets.Item("2021-Q4 (2)")
$newQ1.Name = "2022-Q1"

# ---- Fill the "2022-Q1" sheet with the fund-holdings detail table ----
$q1Header = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($j = 0; $j -lt $q1Header.Count; $j++) {
    $newQ1.Cells.Item(1, 2 + $j).Value = $q1Header[$j]
}

$q1Rows = @(
    @("000849","汇丰晋信双核策略混合A","3.00","85.44","5.92","0.1776","1"),
    @("561550","华泰柏瑞中证500增强策略ETF","9.68","98.93","1.53","0.1481","1"),
    @("012080","易方达中证500指数量化增强型证券投资基金A","6.82","84.83","0.99","0.0675","3"),
    @("010153","中加中证500指数增强A","1.26","94.19","2.41","0.0304","4"),
    @("510200","汇安上证证券ETF","0.74","97.70","3.25","0.0240","8"),
    @("000850","汇丰晋信双核策略混合C","0.39","85.44","5.92","0.0231","1"),
    @("012081","易方达中证500指数量化增强型证券投资基金C","1.57","84.83","0.99","0.0155","3"),
    @("010154","中加中证500指数增强C","0.60","94.19","2.41","0.0145","4"),
    @("165511","信诚中证500指数（LOF）A","2.78","93.31","0.50","0.0139","7"),
    @("510440","大成中证500沪市ETF","0.41","96.76","0.91","0.0037","8"),
    @("006611","人保中证500指数","0.44","92.48","0.59","0.0026","7"),
    @("515550","中融中证500ETF","0.23","91.02","0.47","0.0011","10"),
    @("510570","兴业中证500ETF","0.10","96.12","0.60","0.0006","5"),
    @("013119","信诚中证500指数（LOF）C","0.12","93.31","0.50","0.0006","7")
)

$r = 2
foreach ($row in $q1Rows) {
    $newQ1.Cells.Item($r, 1).Value = $r - 2

    for ($j = 0; $j -lt 6; $j++) {
        $cell = $newQ1.Cells.Item($r, 2 + $j)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$j]
        $cell.ClearFormats()
    }

    $rankCell = $newQ1.Cells.Item($r, 8)
    $rankCell.Value = [int]$row[6]
    $rankCell.ClearFormats()

    $r++
}

# The copied "2021-Q4" sheet had 16 data rows; "2022-Q1" only needs 14,
# so delete the two leftover rows at the bottom (rows 16-17).
$newQ1.Rows.Item(16).Delete()
$newQ1.Rows.Item(16).Delete()

# ---- Rebuild "总计": insert a new first data row for 2022-Q1 ----
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 14
$totalSheet.Range("D2").Value = 0.52

# Re-number the index column (A) for the rows that shifted down.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4

# Keep the original first sheet as the active tab, matching the
# workbook's prior activation state.
$wb.Worksheets.Item("2021-Q1").Activate()
